$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at 295-296, pushing the existing data (old rows 295-333)
# down to rows 297-335.
$ws.Rows("295:296").Insert()

# New row 295: Vega Monumental Concepción, Plátano, Maduro, fecha 44476
$ws.Cells.Item(295, 1).Value = 11
$ws.Cells.Item(295, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(295, 3).Value = "Bíobío"
$ws.Cells.Item(295, 4).Value = 44476
$ws.Cells.Item(295, 5).Value = 8
$ws.Cells.Item(295, 6).Value = "Fruta"
$ws.Cells.Item(295, 7).Value = 100108
$ws.Cells.Item(295, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(295, 9).Value = 100108006
$ws.Cells.Item(295, 10).Value = "Plátano"
$ws.Cells.Item(295, 11).Value = "Sin especificar"
$ws.Cells.Item(295, 12).Value = "Maduro"
$ws.Cells.Item(295, 13).Value = 100
$ws.Cells.Item(295, 14).Value = 18000
$ws.Cells.Item(295, 15).Value = 18000
$ws.Cells.Item(295, 16).Value = 18000
$ws.Cells.Item(295, 17).Value = "$/caja 20 kilos"
$ws.Cells.Item(295, 18).Value = "Ecuador"
$ws.Cells.Item(295, 19).Value = 900
$ws.Cells.Item(295, 20).Value = 20

# New row 296: Vega Monumental Concepción, Plátano, Pintón, fecha 44476
$ws.Cells.Item(296, 1).Value = 11
$ws.Cells.Item(296, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(296, 3).Value = "Bíobío"
$ws.Cells.Item(296, 4).Value = 44476
$ws.Cells.Item(296, 5).Value = 8
$ws.Cells.Item(296, 6).Value = "Fruta"
$ws.Cells.Item(296, 7).Value = 100108
$ws.Cells.Item(296, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(296, 9).Value = 100108006
$ws.Cells.Item(296, 10).Value = "Plátano"
$ws.Cells.Item(296, 11).Value = "Sin especificar"
$ws.Cells.Item(296, 12).Value = "Pintón"
$ws.Cells.Item(296, 13).Value = 300
$ws.Cells.Item(296, 14).Value = 20000
$ws.Cells.Item(296, 15).Value = 20000
$ws.Cells.Item(296, 16).Value = 20000
$ws.Cells.Item(296, 17).Value = "$/caja 20 kilos"
$ws.Cells.Item(296, 18).Value = "Ecuador"
$ws.Cells.Item(296, 19).Value = 1000
$ws.Cells.Item(296, 20).Value = 20
